# Updated RAD Test Cases for Year dropdown
# Records a Katalon-style test execution result (Result / Date) into the
# first two columns of row 2, which previously had no values there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Not Executed"
$ws.Range("B2").Value = "Mon Jan 27 15:34:37 EST 2025"

# These two new cells should keep the sheet's default (unstyled) formatting
# rather than inheriting column A/B's configured style.
$ws.Range("A2:B2").Style = "Normal"
